# Re did the headers overall
# Append new log entries to the "Logs" worksheet (rows 803-812 and 817-818),
# matching the data added in the source commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Logs")

# Row 803
$ws.Cells.Item(803,1).Value = "Pickup PC"
$ws.Cells.Item(803,2).Value = 42669
$ws.Cells.Item(803,3).Value = "1730"
$ws.Cells.Item(803,4).Value = "ATK"
$ws.Cells.Item(803,5).Value = "005"
$ws.Cells.Item(803,6).Value = "Return to DB 0003 !!!!!"

# Row 804
$ws.Cells.Item(804,1).Value = "Pickup Projector"
$ws.Cells.Item(804,2).Value = 42669
$ws.Cells.Item(804,3).Value = "1730"
$ws.Cells.Item(804,4).Value = "ATK"
$ws.Cells.Item(804,5).Value = "005"
$ws.Cells.Item(804,6).Value = "Return to ATK 003C"

# Row 805
$ws.Cells.Item(805,1).Value = "Demo"
$ws.Cells.Item(805,2).Value = 42669
$ws.Cells.Item(805,3).Value = "1900"
$ws.Cells.Item(805,4).Value = "SSB"
$ws.Cells.Item(805,5).Value = "W133"

# Row 806 (taller row, 45pt). The "Special Instructions" text for this row
# (F806) is filled in further below -- after rows 807/808 -- so that new
# shared-string entries are introduced in the same order as the source
# workbook (E806's room number comes first, then F807/F808's comments,
# and only then F806's long comment).
$ws.Cells.Item(806,1).Value = "Setup Mic"
$ws.Cells.Item(806,2).Value = 42669
$ws.Cells.Item(806,3).Value = "1730"
$ws.Cells.Item(806,4).Value = "OSG"
$ws.Cells.Item(806,5).Value = "1017"
$ws.Rows.Item(806).RowHeight = 45

# Row 807
$ws.Cells.Item(807,1).Value = "Pickup Mic"
$ws.Cells.Item(807,2).Value = 42669
$ws.Cells.Item(807,3).Value = "2100"
$ws.Cells.Item(807,4).Value = "OSG"
$ws.Cells.Item(807,5).Value = "1017"
$ws.Cells.Item(807,6).Value = "Return mic, stand, cable , speaker and tripod to OSG 1014L"

# Row 808
$ws.Cells.Item(808,1).Value = "Pickup Skype Kit"
$ws.Cells.Item(808,2).Value = 42669
$ws.Cells.Item(808,3).Value = "1630"
$ws.Cells.Item(808,4).Value = "OSG"
$ws.Cells.Item(808,5).Value = "2010"
$ws.Cells.Item(808,6).Value = "Return Skype kit to OSG 1014L"

# Now fill in F806 (its shared string is introduced after F807/F808's).
$ws.Cells.Item(806,6).Value = "Powered JBL speaker and tripod from HNES 003; mic , stand and cable from OSG 1014L.      ALSO PROVIDE HDMI CABLE FROM OSG 1014L"

# Row 809
$ws.Cells.Item(809,1).Value = "SCLD Student Event"
$ws.Cells.Item(809,2).Value = 42669
$ws.Cells.Item(809,3).Value = "1800"
$ws.Cells.Item(809,4).Value = "R"
$ws.Cells.Item(809,5).Value = "S101"
$ws.Cells.Item(809,6).Value = "INC000000738947"

# Row 810
$ws.Cells.Item(810,1).Value = "SCLD Student Logout"
$ws.Cells.Item(810,2).Value = 42669
$ws.Cells.Item(810,3).Value = "2000"
$ws.Cells.Item(810,4).Value = "R"
$ws.Cells.Item(810,5).Value = "S101"
$ws.Cells.Item(810,6).Value = "INC000000738947"

# Row 811
$ws.Cells.Item(811,1).Value = "SCLD Student Event"
$ws.Cells.Item(811,2).Value = 42669
$ws.Cells.Item(811,3).Value = "1800"
$ws.Cells.Item(811,4).Value = "ATK"
$ws.Cells.Item(811,5).Value = "004"
$ws.Cells.Item(811,6).Value = "INC000000737819"

# Row 812
$ws.Cells.Item(812,1).Value = "SCLD Student Logout"
$ws.Cells.Item(812,2).Value = 42669
$ws.Cells.Item(812,3).Value = "2150"
$ws.Cells.Item(812,4).Value = "ATK"
$ws.Cells.Item(812,5).Value = "004"
$ws.Cells.Item(812,6).Value = "INC000000737819"

# Row 817 (note the gap between 812 and 817 - rows 813-816 remain empty)
$ws.Cells.Item(817,1).Value = "Demo"
$ws.Cells.Item(817,2).Value = 42670
$ws.Cells.Item(817,3).Value = "1800"
$ws.Cells.Item(817,4).Value = "ACE"
$ws.Cells.Item(817,5).Value = "004"

# Row 818
$ws.Cells.Item(818,1).Value = "Demo"
$ws.Cells.Item(818,2).Value = 42670
$ws.Cells.Item(818,3).Value = "1900"
$ws.Cells.Item(818,4).Value = "SSB"
$ws.Cells.Item(818,5).Value = "N108"

# Update the view to reflect scrolling to the bottom of the log and the
# active cell on the next blank row, as in the edited workbook.
$ws.Range("F822").Select()
